$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores Price/Volume figures as literal text (e.g. "60.316.43",
# "1.00", "  -2.89%  ") rather than numbers. Excel's normal type-inference would
# coerce numeric-looking strings like "1.00" or "140.80" into the number 1 / 140.8
# and drop the formatting, so each touched cell is temporarily forced to the Text
# format before its value is written, then the format override is cleared again
# so the cell keeps its original (unstyled) appearance.
$updates = [ordered]@{
    "D2" = "60.316.43"
    "E2" = "  -2.89%  "
    "D3" = "3.300.06"
    "E3" = "  -3.64%  "
    "D4" = "1.00"
    "D5" = "557.73"
    "E5" = "  -3.84%  "
    "D6" = "140.80"
    "E6" = "  -9.02%  "
    "E7" = "  +0.01%  "
    "D8" = "3.299.88"
    "E8" = "  -3.66%  "
    "D9" = "0.467"
    "E9" = "  -3.71%  "
    "D10" = "7.93"
    "E10" = "  -1.92%  "
    "E11" = "  -5.39%  "
    "D12" = "0.407"
    "E12" = "  -2.79%  "
    "D13" = "3.864.50"
    "E13" = "  -3.71%  "
    "E14" = "  -0.32%  "
    "D15" = "26.55"
    "E15" = "  -7.74%  "
    "D16" = "3.295.55"
    "E16" = "  -3.87%  "
    "D17" = "0.0000164"
    "E17" = "  -5.13%  "
    "D18" = "60.286.41"
    "E18" = "  -3.01%  "
    "D19" = "6.09"
    "E19" = "  -6.96%  "
    "D20" = "13.65"
    "E20" = "  -5.68%  "
    "D21" = "8.52"
    "E21" = "  -5.40%  "
    "D22" = "373.78"
    "E22" = "  -2.65%  "
    "E23" = "  -0.08%  "
    "D24" = "72.07"
    "E24" = "  -5.23%  "
    "D25" = "0.532"
    "E25" = "  -7.11%  "
    "D26" = "3.430.82"
    "E26" = "  -3.71%  "
    "E27" = "  -9.71%  "
    "E28" = "  -1.51%  "
    "D29" = "1.00"
    "E29" = "  +0.16%  "
    "D30" = "7.03"
    "E30" = "  -8.60%  "
    "E31" = "  -0.06%  "
    "E32" = "  -5.11%  "
    "D33" = "7.27"
    "E33" = "  -7.95%  "
    "D34" = "22.54"
    "E34" = "  -3.37%  "
    "D35" = "1.26"
    "E35" = "  -5.14%  "
    "D36" = "165.37"
    "E36" = "  -1.96%  "
    "D37" = "5.02"
    "E37" = "  -9.69%  "
    "E38" = "  -5.19%  "
    "D39" = "6.62"
    "E39" = "  -5.35%  "
    "D40" = "3.332.41"
    "E40" = "  -3.75%  "
    "D41" = "0.0723"
    "E41" = "  -7.75%  "
    "D42" = "25.15"
    "E42" = "  -18.81%  "
    "D43" = "41.80"
    "E43" = "  -2.32%  "
    "D44" = "0.748"
    "E44" = "  -4.27%  "
    "E45" = "  -4.32%  "
    "D46" = "4.09"
    "E46" = "  -7.46%  "
    "D47" = "1.57"
    "E47" = "  -6.69%  "
    "E48" = "  -0.05%  "
    "D49" = "2.318.27"
    "E49" = "  -9.33%  "
    "D50" = "6.37"
    "E50" = "  -6.65%  "
    "D51" = "21.45"
    "E51" = "  -8.48%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).ClearFormats()
}
